# Auto-generated edit script: applies the cryptos list price/volume refresh
# described in the commit "Updated cryptos list on Fri Feb 23 05:11:04 UTC 2024
# with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '51.148.70'
$ws.Cells.Item(2, 5).Value = '  -0.55%  '

$ws.Cells.Item(3, 4).Value = '2.948.95'
$ws.Cells.Item(3, 5).Value = '  +0.71%  '

$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$ws.Cells.Item(5, 4).Value = '''378.45'
$ws.Cells.Item(5, 5).Value = '  +0.78%  '

$ws.Cells.Item(6, 4).Value = '''102.19'
$ws.Cells.Item(6, 5).Value = '  -1.15%  '

$ws.Cells.Item(7, 5).Value = '  -0.44%  '

$ws.Cells.Item(8, 5).Value = '  +0.02%  '

$ws.Cells.Item(9, 4).Value = '''0.588'
$ws.Cells.Item(9, 5).Value = '  +0.53%  '

$ws.Cells.Item(10, 5).Value = '  -0.84%  '

$ws.Cells.Item(11, 5).Value = '  -0.12%  '

$ws.Cells.Item(12, 4).Value = '''0.0839'
$ws.Cells.Item(12, 5).Value = '  +0.22%  '

$ws.Cells.Item(13, 4).Value = '3.416.02'
$ws.Cells.Item(13, 5).Value = '  +0.84%  '

$ws.Cells.Item(14, 4).Value = '''18.01'
$ws.Cells.Item(14, 5).Value = '  -1.68%  '

$ws.Cells.Item(15, 4).Value = '''7.41'
$ws.Cells.Item(15, 5).Value = '  +0.47%  '

$ws.Cells.Item(16, 4).Value = '2.947.56'
$ws.Cells.Item(16, 5).Value = '  +0.83%  '

$ws.Cells.Item(17, 4).Value = '''0.986'
$ws.Cells.Item(17, 5).Value = '  +5.06%  '

$ws.Cells.Item(18, 4).Value = '51.094.91'

$ws.Cells.Item(19, 4).Value = '''3.20'
$ws.Cells.Item(19, 5).Value = '  -5.76%  '

$ws.Cells.Item(20, 4).Value = '''7.21'
$ws.Cells.Item(20, 5).Value = '  -1.31%  '

$ws.Cells.Item(21, 4).Value = '''12.55'
$ws.Cells.Item(21, 5).Value = '  -3.18%  '

$ws.Cells.Item(22, 5).Value = '  +0.59%  '

$ws.Cells.Item(23, 4).Value = '''68.43'
$ws.Cells.Item(23, 5).Value = '  +0.26%  '

$ws.Cells.Item(24, 4).Value = '''261.80'
$ws.Cells.Item(24, 5).Value = '  +0.19%  '

$ws.Cells.Item(25, 5).Value = '  +2.52%  '

$ws.Cells.Item(26, 4).Value = '''8.41'
$ws.Cells.Item(26, 5).Value = '  +14.23%  '

$ws.Cells.Item(27, 4).Value = '''7.59'
$ws.Cells.Item(27, 5).Value = '  +7.25%  '

$ws.Cells.Item(28, 5).Value = '  +0.21%  '

$ws.Cells.Item(29, 5).Value = '  -0.43%  '

$ws.Cells.Item(30, 5).Value = '  +11.83%  '

$ws.Cells.Item(31, 5).Value = '  -0.08%  '

$ws.Cells.Item(32, 4).Value = '''25.68'
$ws.Cells.Item(32, 5).Value = '  -0.39%  '

$ws.Cells.Item(33, 5).Value = '  -0.13%  '

$ws.Cells.Item(34, 4).Value = '''33.82'
$ws.Cells.Item(34, 5).Value = '  -0.49%  '

$ws.Cells.Item(35, 4).Value = '''50.36'
$ws.Cells.Item(35, 5).Value = '  -2.87%  '

$ws.Cells.Item(36, 5).Value = '  -2.04%  '

$ws.Cells.Item(37, 4).Value = '''0.0449'
$ws.Cells.Item(37, 5).Value = '  +5.45%  '

$ws.Cells.Item(38, 5).Value = '  -0.01%  '

$ws.Cells.Item(39, 5).Value = '  -1.65%  '

$ws.Cells.Item(40, 4).Value = '''16.93'
$ws.Cells.Item(40, 5).Value = '  -0.06%  '

$ws.Cells.Item(41, 4).Value = '''2.56'
$ws.Cells.Item(41, 5).Value = '  -1.40%  '

$ws.Cells.Item(42, 5).Value = '  +0.76%  '

$ws.Cells.Item(43, 5).Value = '  -2.22%  '

$ws.Cells.Item(44, 4).Value = '''121.66'
$ws.Cells.Item(44, 5).Value = '  -2.15%  '

$ws.Cells.Item(45, 4).Value = '''21.06'
$ws.Cells.Item(45, 5).Value = '  -3.50%  '

$ws.Cells.Item(46, 5).Value = '  -0.17%  '

$ws.Cells.Item(47, 5).Value = '  +1.73%  '

$ws.Cells.Item(48, 5).Value = '  +2.03%  '

$ws.Cells.Item(49, 2).Value = 'Maker'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(49, 4).Value = '2.000.45'
$ws.Cells.Item(49, 5).Value = '  -0.92%  '

$ws.Cells.Item(50, 2).Value = 'NEARProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(50, 4).Value = '''3.22'
$ws.Cells.Item(50, 5).Value = '  +1.34%  '

$ws.Cells.Item(51, 4).Value = '''0.0336'
$ws.Cells.Item(51, 5).Value = '  +4.82%  '
